# The upstream export re-sorted rows 19-30 of the sheet (a new batch of
# validation records pushed each record's row down/up, row 29 happening to
# land back on itself). Per record, only the "identity" columns travel -
# Id (A), Taxonsorteringsordning (B), Rodlistade (D), TaxonId (E),
# Artnamn (F), Vetenskapligt namn (G), Auktor (H), Kon (L, only ever
# empty), and the Ost/Nord coordinates (Q/R). Every other column in the
# row (date, observer, municipality, ...) already repeats identically
# across this record block, so it is left completely untouched rather
# than being rewritten - that avoids Excel's COM layer "helpfully"
# reinterpreting the plain-text dates in Y/AA as real dates, and avoids
# disturbing the handful of genuinely-empty cells elsewhere in the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row -> old row that now supplies its identity/location data.
$mapping = @{
    19 = 20
    20 = 30
    21 = 24
    22 = 23
    23 = 26
    24 = 19
    25 = 22
    26 = 21
    27 = 28
    28 = 25
    29 = 29
    30 = 27
}

$movingCols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Snapshot every row's moving columns first, so later writes never clobber
# a row we still need to read from (row 29 maps to itself, and several
# rows swap with each other).
$snapshots = @{}
foreach ($r in $mapping.Keys) {
    $rowVals = @{}
    foreach ($col in $movingCols) {
        $rowVals[$col] = $ws.Range($col + $r).Value2
    }
    $snapshots[$r] = $rowVals
}

# Column L ("Kon") is an empty placeholder in every row in this block -
# only its bare presence/absence as an empty cell toggles between rows.
# Record which source rows had it so the destination can be cleared to
# match (Excel COM has no way to write a "present but empty" cell - an
# empty assignment always clears the cell outright, which is exactly the
# behaviour wanted for the destinations below).
$lPresent = @{
    19 = $true;  20 = $false; 21 = $true;  22 = $false
    23 = $true;  24 = $false; 25 = $false; 26 = $true
    27 = $true;  28 = $false; 29 = $false; 30 = $false
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $rowVals = $snapshots[$oldRow]

    foreach ($col in $movingCols) {
        $ws.Range($col + $newRow).Value2 = $rowVals[$col]
    }

    if (-not $lPresent[$oldRow]) {
        $ws.Range("L" + $newRow).Value2 = ""
    }
}
